$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.06924433320117
$ws.Range("D2").Value = 4.485380041099008
$ws.Range("E2").Value = 18.29363633744745
$ws.Range("F2").Value = 21.51286486881103
$ws.Range("G2").Value = 23.98580428202334
$ws.Range("H2").Value = 13.10355686862735
$ws.Range("I2").Value = 25.35969145240661
$ws.Range("K2").Value = 10.6471623488218
$ws.Range("L2").Value = 8.731502971933187
$ws.Range("M2").Value = 13.89023653698495
$ws.Range("N2").Value = 19.37426800582106
$ws.Range("O2").Value = 19.30405114390045
$ws.Range("B3").Value = 12.96024306971328
$ws.Range("D3").Value = 4.401481974132078
$ws.Range("E3").Value = 18.34765482451671
$ws.Range("F3").Value = 21.51568994599208
$ws.Range("G3").Value = 23.98416585303979
$ws.Range("H3").Value = 13.13832966020087
$ws.Range("I3").Value = 25.46125114090707
$ws.Range("K3").Value = 10.39280996146833
$ws.Range("L3").Value = 8.71640375954947
$ws.Range("M3").Value = 13.86540324609726
$ws.Range("N3").Value = 19.42594937714618
$ws.Range("O3").Value = 19.35137978003875
$ws.Range("B4").Value = 12.89496880521022
$ws.Range("D4").Value = 4.3484008653672
$ws.Range("E4").Value = 18.38272950782679
$ws.Range("F4").Value = 21.52288157414488
$ws.Range("G4").Value = 23.99070893787583
$ws.Range("H4").Value = 13.16156827907988
$ws.Range("I4").Value = 25.5273002452715
$ws.Range("K4").Value = 10.23181876199396
$ws.Range("L4").Value = 8.708352780309589
$ws.Range("M4").Value = 13.85192404418803
$ws.Range("N4").Value = 19.45929581293369
$ws.Range("O4").Value = 19.38430473405208
$ws.Range("B5").Value = 12.86880968361858
$ws.Range("D5").Value = 4.326390220291992
$ws.Range("E5").Value = 18.3975034773938
$ws.Range("F5").Value = 21.52718537314339
$ws.Range("G5").Value = 23.99527317130263
$ws.Range("H5").Value = 13.17151309837521
$ws.Range("I5").Value = 25.55514534262433
$ws.Range("K5").Value = 10.16505884110655
$ws.Range("L5").Value = 8.705381296979745
$ws.Range("M5").Value = 13.84688009318463
$ws.Range("N5").Value = 19.47329161990058
$ws.Range("O5").Value = 19.39869264064093
$ws.Range("B6").Value = 12.86449331346763
$ws.Range("D6").Value = 4.322712875307767
$ws.Range("E6").Value = 18.39998575511785
$ws.Range("F6").Value = 21.52798296181601
$ws.Range("G6").Value = 23.99614564389354
$ws.Range("H6").Value = 13.17319311451313
$ws.Range("I6").Value = 25.5598251842191
$ws.Range("K6").Value = 10.15390539862256
$ws.Range("L6").Value = 8.704906640297544
$ws.Range("M6").Value = 13.84606978143823
$ws.Range("N6").Value = 19.47564021735171
$ws.Range("O6").Value = 19.4011403341818
$ws.Range("B7").Value = 12.89461419737818
$ws.Range("D7").Value = 4.348105539729014
$ws.Range("E7").Value = 18.38292680660327
$ws.Range("F7").Value = 21.52293405621753
$ws.Range("G7").Value = 23.99076280997733
$ws.Range("H7").Value = 13.1617004753798
$ws.Range("I7").Value = 25.52767200838219
$ws.Range("K7").Value = 10.23092301193638
$ws.Range("L7").Value = 8.708311450012204
$ws.Range("M7").Value = 13.85185419661028
$ws.Range("N7").Value = 19.45948291651339
$ws.Range("O7").Value = 19.38449484561322
$ws.Range("B8").Value = 13.03133285618357
$ws.Range("D8").Value = 4.456785743106589
$ws.Range("E8").Value = 18.31186680859607
$ws.Range("F8").Value = 21.51270716014581
$ws.Range("G8").Value = 23.98367274953587
$ws.Range("H8").Value = 13.11515476626303
$ws.Range("I8").Value = 25.39394400359402
$ws.Range("K8").Value = 10.56049447365032
$ws.Range("L8").Value = 8.726045154695306
$ws.Range("M8").Value = 13.88130956598467
$ws.Range("N8").Value = 19.39175341507602
$ws.Range("O8").Value = 19.31956722202699
$ws.Range("B9").Value = 13.3114224389465
$ws.Range("D9").Value = 4.656874079574178
$ws.Range("E9").Value = 18.18759618227186
$ws.Range("F9").Value = 21.53587514517894
$ws.Range("G9").Value = 24.02961823933182
$ws.Range("H9").Value = 13.03885247344874
$ws.Range("I9").Value = 25.16092622732621
$ws.Range("K9").Value = 11.16625612659289
$ws.Range("L9").Value = 8.770380412890191
$ws.Range("M9").Value = 13.95290867793949
$ws.Range("N9").Value = 19.2716932574986
$ws.Range("O9").Value = 19.22295863257616
$ws.Range("B10").Value = 13.52296049736657
$ws.Range("D10").Value = 4.795189347689069
$ws.Range("E10").Value = 18.1054105883442
$ws.Range("F10").Value = 21.57909281739151
$ws.Range("G10").Value = 24.09969372003662
$ws.Range("H10").Value = 12.99191222888087
$ws.Range("I10").Value = 25.00745266921532
$ws.Range("K10").Value = 11.58379523528498
$ws.Range("L10").Value = 8.808608961329664
$ws.Range("M10").Value = 14.01367192020838
$ws.Range("N10").Value = 19.19119420456903
$ws.Range("O10").Value = 19.17076253089757
$ws.Range("B11").Value = 13.62010802239498
$ws.Range("D11").Value = 4.856083541605948
$ws.Range("E11").Value = 18.06998564500184
$ws.Range("F11").Value = 21.60439080353468
$ws.Range("G11").Value = 24.13938999321795
$ws.Range("H11").Value = 12.97253586180224
$ws.Range("I11").Value = 24.94146356133702
$ws.Range("K11").Value = 11.76725159796561
$ws.Range("L11").Value = 8.827187595347601
$ws.Range("M11").Value = 14.04302230893867
$ws.Range("N11").Value = 19.1562326645736
$ws.Range("O11").Value = 19.15110492015229
$ws.Range("B12").Value = 13.65699879069764
$ws.Range("D12").Value = 4.878840326636509
$ws.Range("E12").Value = 18.05685202670864
$ws.Range("F12").Value = 21.61477540980284
$ws.Range("G12").Value = 24.15553773066539
$ws.Range("H12").Value = 12.96548265045547
$ws.Range("I12").Value = 24.91702417847883
$ws.Range("K12").Value = 11.83575077247964
$ws.Range("L12").Value = 8.834389968619906
$ws.Range("M12").Value = 14.05437641045972
$ws.Range("N12").Value = 19.14323098406062
$ws.Range("O12").Value = 19.1442492278876
$ws.Range("B13").Value = 13.64904957662977
$ws.Range("D13").Value = 4.873952870354127
$ws.Range("E13").Value = 18.0596681037258
$ws.Range("F13").Value = 21.61250320914807
$ws.Range("G13").Value = 24.15201056461905
$ws.Range("H13").Value = 12.96698904750698
$ws.Range("I13").Value = 24.92226322435664
$ws.Range("K13").Value = 11.82104205850964
$ws.Range("L13").Value = 8.832831441927606
$ws.Range("M13").Value = 14.05192053666958
$ws.Range("N13").Value = 19.14602058018018
$ws.Range("O13").Value = 19.14569955569397
$ws.Range("B14").Value = 13.62314114671401
$ws.Range("D14").Value = 4.857961880299221
$ws.Range("E14").Value = 18.06889950785703
$ws.Range("F14").Value = 21.60522905901447
$ws.Range("G14").Value = 24.14069617415049
$ws.Range("H14").Value = 12.97194989466812
$ws.Range("I14").Value = 24.93944191636998
$ws.Range("K14").Value = 11.77290675485394
$ws.Range("L14").Value = 8.827776813572646
$ws.Range("M14").Value = 14.04395164798608
$ws.Range("N14").Value = 19.15515825386256
$ws.Range("O14").Value = 19.15052910687655
$ws.Range("B15").Value = 13.60728405763978
$ws.Range("D15").Value = 4.848127221136493
$ws.Range("E15").Value = 18.07459057802586
$ws.Range("F15").Value = 21.60087805011888
$ws.Range("G15").Value = 24.13391078435856
$ws.Range("H15").Value = 12.97502556489445
$ws.Range("I15").Value = 24.95003586124103
$ws.Range("K15").Value = 11.74329483086441
$ws.Range("L15").Value = 8.824702345811382
$ws.Range("M15").Value = 14.03910151377501
$ws.Range("N15").Value = 19.16078624515813
$ws.Range("O15").Value = 19.15356396340839
$ws.Range("B16").Value = 13.51662767744892
$ws.Range("D16").Value = 4.791168025644978
$ws.Range("E16").Value = 18.10776507596857
$ws.Range("F16").Value = 21.57755249971608
$ws.Range("G16").Value = 24.09725610040439
$ws.Range("H16").Value = 12.99321829470111
$ws.Range("I16").Value = 25.0118421287477
$ws.Range("K16").Value = 11.57167200724109
$ws.Range("L16").Value = 8.807418392111909
$ws.Range("M16").Value = 14.011787683947
$ws.Range("N16").Value = 19.19351230656451
$ws.Range("O16").Value = 19.17212948812788
$ws.Range("B17").Value = 13.46122744125111
$ws.Range("D17").Value = 4.755698272949734
$ws.Range("E17").Value = 18.12861824206075
$ws.Range("F17").Value = 21.5646831344497
$ws.Range("G17").Value = 24.07676651054246
$ws.Range("H17").Value = 13.0048852393293
$ws.Range("I17").Value = 25.05073766189257
$ws.Range("K17").Value = 11.46469707402173
$ws.Range("L17").Value = 8.797116927086817
$ws.Range("M17").Value = 13.99546503622776
$ws.Range("N17").Value = 19.21401266039265
$ws.Range("O17").Value = 19.18456598392678
$ws.Range("B18").Value = 13.42945080676015
$ws.Range("D18").Value = 4.735106920263192
$ws.Range("E18").Value = 18.14079713782543
$ws.Range("F18").Value = 21.55781210908615
$ws.Range("G18").Value = 24.06571798730959
$ws.Range("H18").Value = 13.01178184680111
$ws.Range("I18").Value = 25.07346959806835
$ws.Range("K18").Value = 11.40255971443447
$ws.Range("L18").Value = 8.791303883269865
$ws.Range("M18").Value = 13.98623790337028
$ws.Range("N18").Value = 19.22596002281893
$ws.Range("O18").Value = 19.19210376055194
$ws.Range("B19").Value = 13.41870779880911
$ws.Range("D19").Value = 4.728102747151082
$ws.Range("E19").Value = 18.14495245832355
$ws.Range("F19").Value = 21.55557707973944
$ws.Range("G19").Value = 24.06210388748481
$ws.Range("H19").Value = 13.01414888314429
$ws.Range("I19").Value = 25.08122814691695
$ws.Range("K19").Value = 11.38141786841591
$ws.Range("L19").Value = 8.789355051109963
$ws.Range("M19").Value = 13.98314162335264
$ws.Range("N19").Value = 19.23003203277726
$ws.Range("O19").Value = 19.19472196326717
$ws.Range("B20").Value = 13.46711598061332
$ws.Range("D20").Value = 4.759493845549601
$ws.Range("E20").Value = 18.12637927819162
$ws.Range("F20").Value = 21.56599817432958
$ws.Range("G20").Value = 24.07887149020604
$ws.Range("H20").Value = 13.00362401419437
$ws.Range("I20").Value = 25.04655988857185
$ws.Range("K20").Value = 11.47614798155917
$ws.Range("L20").Value = 8.798201962252634
$ws.Range("M20").Value = 13.99718596665248
$ws.Range("N20").Value = 19.21181421345059
$ws.Range("O20").Value = 19.18320228498548
$ws.Range("B21").Value = 13.63074851192543
$ws.Range("D21").Value = 4.86266711734424
$ws.Range("E21").Value = 18.0661804016908
$ws.Range("F21").Value = 21.6073438627815
$ws.Range("G21").Value = 24.14398928458731
$ws.Range("H21").Value = 12.97048506133114
$ws.Range("I21").Value = 24.93438121936144
$ws.Range("K21").Value = 11.78707193002728
$ws.Range("L21").Value = 8.829256978164036
$ws.Range("M21").Value = 14.04628584603164
$ws.Range("N21").Value = 19.15246785918456
$ws.Range("O21").Value = 19.14909458296662
$ws.Range("B22").Value = 13.73827952335727
$ws.Range("D22").Value = 4.928328828717062
$ws.Range("E22").Value = 18.0284745866385
$ws.Range("F22").Value = 21.6390533065356
$ws.Range("G22").Value = 24.19304504053244
$ws.Range("H22").Value = 12.95048331315741
$ws.Range("I22").Value = 24.8642670301176
$ws.Range("K22").Value = 11.98459868305366
$ws.Range("L22").Value = 8.850525013102057
$ws.Range("M22").Value = 14.07976993616882
$ws.Range("N22").Value = 19.11506548277093
$ws.Range("O22").Value = 19.13023194632827
$ws.Range("B23").Value = 13.68084393471173
$ws.Range("D23").Value = 4.893449200570663
$ws.Range("E23").Value = 18.04844939049032
$ws.Range("F23").Value = 21.62170265432785
$ws.Range("G23").Value = 24.16627182716329
$ws.Range("H23").Value = 12.96100708868446
$ws.Range("I23").Value = 24.90139572612324
$ws.Range("K23").Value = 11.8797067587743
$ws.Range("L23").Value = 8.839086252250993
$ws.Range("M23").Value = 14.06177329181132
$ws.Range("N23").Value = 19.13490150422869
$ws.Range("O23").Value = 19.13998542386157
$ws.Range("B24").Value = 13.46445354020691
$ws.Range("D24").Value = 4.757778487161334
$ws.Range("E24").Value = 18.12739092144837
$ws.Range("F24").Value = 21.56540200049286
$ws.Range("G24").Value = 24.07791755060019
$ws.Range("H24").Value = 13.00419362495891
$ws.Range("I24").Value = 25.04844750612837
$ws.Range("K24").Value = 11.47097300440408
$ws.Range("L24").Value = 8.797711076781908
$ws.Range("M24").Value = 13.99640744460041
$ws.Range("N24").Value = 19.21280762838346
$ws.Range("O24").Value = 19.18381760514122
$ws.Range("B25").Value = 13.23452619625492
$ws.Range("D25").Value = 4.604219996846442
$ws.Range("E25").Value = 18.219608374463
$ws.Range("F25").Value = 21.52499261723161
$ws.Range("G25").Value = 24.0107922411158
$ws.Range("H25").Value = 13.05789214641342
$ws.Range("I25").Value = 25.22084469081848
$ws.Range("K25").Value = 11.00702485261052
$ws.Range("L25").Value = 8.757380160262031
$ws.Range("M25").Value = 13.93208529967463
$ws.Range("N25").Value = 19.30281404077451
$ws.Range("O25").Value = 19.24579903704835
